$d = $word.ActiveDocument
$p = $d.Paragraphs.Item(138)
$rng = $p.Range
$rng.Collapse(1)
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Test Objective:</w:t></w:r></w:p><w:p><w:r><w:t>The objective of this test is to ensure that the forums feature of the website is implemented with the highest possible quality. Identifying and safeguarding the user from issues within the forum is a priority to ensure a good user experience. A user should be able to navigate through the forums and easily post a reply when logged-in, regardless of browser type.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Hardware and software setup:</w:t></w:r></w:p><w:p><w:r><w:t>Hardware for QA testing includes a computer, keyboard, mouse, and monitor. Software includes a web browser, specifically Chrome and Firefox.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Feature to be tested:</w:t></w:r></w:p><w:p><w:r><w:t>The ability for a user to navigate to a forum discussion and post a reply will be tested.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Actual test cases:</w:t></w:r></w:p><w:p><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> </w:t></w:r></w:p><w:tbl><w:tblPr><w:tblStyle w:val="TableGrid"/><w:tblW w:w="9895" w:type="dxa"/><w:tblLook w:val="04A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="0" w:noVBand="1"/></w:tblPr><w:tblGrid><w:gridCol w:w="937"/><w:gridCol w:w="1308"/><w:gridCol w:w="1890"/><w:gridCol w:w="1890"/><w:gridCol w:w="2160"/><w:gridCol w:w="1710"/></w:tblGrid><w:tr><w:tc><w:tcPr><w:tcW w:w="937" w:type="dxa"/></w:tcPr><w:p><w:bookmarkStart w:id="14" w:name="_Hlk87883177"/><w:r><w:t>Test Number</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="1308" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>Test Title</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="1890" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>Test Description</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="1890" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>Test Input</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="2160" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>Expected Correct Output</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="1710" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>Test Results</w:t></w:r></w:p></w:tc></w:tr><w:bookmarkEnd w:id="14"/><w:tr><w:tc><w:tcPr><w:tcW w:w="937" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>1</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="1308" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>Forum Navigation</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="1890" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>The user navigates to a forum post.</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="1890" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>Mouse input will click through appropriate links.</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="2160" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>The web interface should display the appropriate discussion.</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="1710" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>PASS – Chrome</w:t></w:r><w:r><w:br/><w:t>PASS – Firefox</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="937" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>2</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="1308" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>Post Valid Reply</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="1890" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>The user attempts to post a valid reply.</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="1890" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>The post will contain “Hello World, I’m testing input”.</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="2160" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>The forum interface should update with a field dedicated to and with the user’s post.</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="1710" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>PASS – Chrome</w:t></w:r><w:r><w:br/><w:t>PASS – Firefox</w:t></w:r><w:r><w:br/></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="937" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>3</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="1308" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>Post Invalid Reply</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="1890" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>The user attempts to post an invalid reply consisting of 0 characters.</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="1890" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>The post will contain nothing.</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="2160" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>The forum should reject an empty post with a warning prompt alerting the user to enter text.</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="1710" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>PASS – Chrome</w:t></w:r><w:r><w:br/><w:t>PASS – Firefox</w:t></w:r><w:r><w:br/></w:r></w:p></w:tc></w:tr></w:tbl></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rng.InsertXML($xml)
Write-Host "Inserted QA test content"
